$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace the underscore signature line that sits right above the
#    ${trabajador} / ${vendedor} row with ${trabajador_linea}, and
#    relocate the single, document-wide _GoBack bookmark so that it
#    sits right after the newly inserted text (mirrors Word leaving
#    the caret there after the last edit).
# ------------------------------------------------------------------
$underscoreRun = "____________________________________"

$hit = $null
$scan = $d.Content
$scan.Find.Execute($underscoreRun, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($scan.Find.Found) {
    $lookahead = $d.Range($scan.End, [Math]::Min($scan.End + 120, $d.Content.End))
    if ($lookahead.Text -like '*${trabajador}*') {
        $hit = $d.Range($scan.Start, $scan.End)
        break
    }
    $scan.Collapse(0)
    $scan.Find.Execute($underscoreRun, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

if ($hit -ne $null) {
    # Insert the replacement text just before the underscores (rather
    # than assigning .Text directly on the run) so the new run's end
    # position is a "clean" insertion point for the bookmark.
    $insertPoint = $d.Range($hit.Start, $hit.Start)
    $insertPoint.InsertBefore('${trabajador_linea}')

    $afterNewText = $insertPoint.End

    $bmRange = $d.Range($afterNewText, $afterNewText)
    $bmRange.Bookmarks.Add("_GoBack")

    # Now drop the original underscores, which were pushed right after
    # the text we just inserted.
    $oldRun = $d.Range($afterNewText, $afterNewText + $underscoreRun.Length)
    $oldRun.Text = ""
}

# ------------------------------------------------------------------
# 2) Drop the redundant "DNI: " prefix that precedes the
#    ${trabajador_dni} placeholder (the ${vendedor_dni} one is left
#    untouched).
# ------------------------------------------------------------------
$prefix = 'DNI: ${'

$scan2 = $d.Content
$scan2.Find.Execute($prefix, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($scan2.Find.Found) {
    $lookahead2 = $d.Range($scan2.End, [Math]::Min($scan2.End + 30, $d.Content.End))
    if ($lookahead2.Text -like 'trabajador_dni*') {
        $scan2.Text = '${'
        break
    }
    $scan2.Collapse(0)
    $scan2.Find.Execute($prefix, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# ------------------------------------------------------------------
# 3) The old _GoBack bookmark (formerly after "SOBRE LAS FUTURAS
#    NOTIFICACIONES JUDICIALES") was already relocated by the
#    Bookmarks.Add() call above, since a document can only contain one
#    bookmark with a given name -- nothing further to delete here.
# ------------------------------------------------------------------
